$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 501
$ws.Range("I6").Value = 320.6
$ws.Range("J6").Value = 726.5
$ws.Range("K6").Value = 961.8000000000001
$ws.Range("L6").Value = 2179.5
$ws.Range("M6").Value = -849.8000000000001
$ws.Range("N6").Value = -2403.5
$ws.Range("H11").Value = 182
$ws.Range("I11").Value = 182
$ws.Range("K11").Value = 182
$ws.Range("M11").Value = -42
$ws.Range("H54").Value = 22515.285
$ws.Range("J54").Value = 22515.285
$ws.Range("L54").Value = 22515.285
$ws.Range("N54").Value = -23487.285
$ws.Range("H55").Value = 190.9
$ws.Range("I55").Value = 262.33334
$ws.Range("J55").Value = 83.75
$ws.Range("K55").Value = 262.33334
$ws.Range("L55").Value = 83.75
$ws.Range("M55").Value = -48.33334000000002
$ws.Range("N55").Value = -511.75
$ws.Range("H101").Value = 9697.5
$ws.Range("I101").Value = 9995
$ws.Range("J101").Value = 9400
$ws.Range("K101").Value = 29985
$ws.Range("L101").Value = 28200
$ws.Range("M101").Value = -28363
$ws.Range("N101").Value = -31444
$ws.Range("H138").Value = 1944.91
$ws.Range("I138").Value = 2108.3333
$ws.Range("J138").Value = 1922.625
$ws.Range("K138").Value = 6324.999899999999
$ws.Range("L138").Value = 5767.875
$ws.Range("M138").Value = -1184.999899999999
$ws.Range("N138").Value = -16047.875

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H13").Value = 9074.833000000001
$ws.Range("J13").Value = 22499.5
$ws.Range("L13").Value = 22499.5
$ws.Range("N13").Value = -22787.5
$ws.Range("H61").Value = 2182.3809
$ws.Range("I61").Value = 1923.6666
$ws.Range("K61").Value = 1923.6666
$ws.Range("M61").Value = -1711.6666
$ws.Range("H102").Value = 8400.25
$ws.Range("I102").Value = 8400.25
$ws.Range("K102").Value = 8400.25
$ws.Range("M102").Value = -6778.25
$ws.Range("H132").Value = 3841.2766
$ws.Range("I132").Value = 3798.0244
$ws.Range("J132").Value = 4136.8335
$ws.Range("K132").Value = 11394.0732
$ws.Range("L132").Value = 12410.5005
$ws.Range("M132").Value = -8864.073199999999
$ws.Range("N132").Value = -17470.5005
$ws.Range("H136").Value = 2182.3809
$ws.Range("I136").Value = 1923.6666
$ws.Range("K136").Value = 5770.9998
$ws.Range("M136").Value = -3220.9998

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 2771.75
$ws.Range("I22").Value = 2544
$ws.Range("K22").Value = 2544
$ws.Range("M22").Value = -2371
$ws.Range("H94").Value = 533.3461
$ws.Range("I94").Value = 272.3889
$ws.Range("K94").Value = 272.3889
$ws.Range("M94").Value = 178.6111

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 582.75
$ws.Range("I16").Value = 582.75
$ws.Range("K16").Value = 582.75
$ws.Range("M16").Value = -295.75
$ws.Range("H31").Value = 38632.5
$ws.Range("I31").Value = 1575.9445
$ws.Range("K31").Value = 1575.9445
$ws.Range("M31").Value = -1280.9445
$ws.Range("H34").Value = 38632.5
$ws.Range("I34").Value = 1575.9445
$ws.Range("K34").Value = 1575.9445
$ws.Range("M34").Value = -1373.9445
$ws.Range("H113").Value = 582.75
$ws.Range("I113").Value = 582.75
$ws.Range("K113").Value = 582.75
$ws.Range("M113").Value = 1587.25
$ws.Range("H141").Value = 284380
$ws.Range("J141").Value = 525360
$ws.Range("L141").Value = 525360
$ws.Range("N141").Value = -535720

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H70").Value = 6140.875
$ws.Range("J70").Value = 6140.875
$ws.Range("L70").Value = 18422.625
$ws.Range("N70").Value = -19052.625
$ws.Range("H73").Value = 6140.875
$ws.Range("J73").Value = 6140.875
$ws.Range("L73").Value = 18422.625
$ws.Range("N73").Value = -20606.625
$ws.Range("H75").Value = 1933.9
$ws.Range("I75").Value = 1799
$ws.Range("J75").Value = 1991.7142
$ws.Range("K75").Value = 5397
$ws.Range("L75").Value = 5975.142599999999
$ws.Range("M75").Value = -4399
$ws.Range("N75").Value = -7971.142599999999
$ws.Range("H78").Value = 1933.9
$ws.Range("I78").Value = 1799
$ws.Range("J78").Value = 1991.7142
$ws.Range("K78").Value = 16191
$ws.Range("L78").Value = 17925.4278
$ws.Range("M78").Value = -11199
$ws.Range("N78").Value = -27909.4278
$ws.Range("H103").Value = 1174.7142
$ws.Range("J103").Value = 654
$ws.Range("L103").Value = 1962
$ws.Range("N103").Value = -3720
$ws.Range("H107").Value = 2087.5
$ws.Range("I107").Value = 0
$ws.Range("J107").Value = 2087.5
$ws.Range("K107").Value = 0
$ws.Range("L107").Value = 6262.5
$ws.Range("M107").ClearContents()
$ws.Range("N107").Value = -10102.5
$ws.Range("H125").Value = 28237.5
$ws.Range("I125").Value = 16000
$ws.Range("K125").Value = 48000
$ws.Range("M125").Value = -43080

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H47").Value = 39936.5
$ws.Range("J47").Value = 39936.5
$ws.Range("L47").Value = 39936.5
$ws.Range("N47").Value = -41072.5
$ws.Range("H52").Value = 30000
$ws.Range("J52").Value = 30000
$ws.Range("L52").Value = 30000
$ws.Range("N52").Value = -30518
$ws.Range("H104").Value = 49500
$ws.Range("J104").Value = 49500
$ws.Range("L104").Value = 49500
$ws.Range("N104").Value = -56488
$ws.Range("H126").Value = 9383.909
$ws.Range("I126").Value = 9999.75
$ws.Range("J126").Value = 9032
$ws.Range("K126").Value = 29999.25
$ws.Range("L126").Value = 27096
$ws.Range("M126").Value = -27529.25
$ws.Range("N126").Value = -32036
$ws.Range("H132").Value = 8189.4165
$ws.Range("I132").Value = 3424.5
$ws.Range("K132").Value = 10273.5
$ws.Range("M132").Value = -7743.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H14").Value = 10000
$ws.Range("J14").Value = 10000
$ws.Range("L14").Value = 10000
$ws.Range("N14").Value = -10344
$ws.Range("H16").Value = 813.2222
$ws.Range("I16").Value = 813.2222
$ws.Range("K16").Value = 813.2222
$ws.Range("M16").Value = -643.2222
$ws.Range("H18").Value = 52000010
$ws.Range("J18").Value = 52000010
$ws.Range("L18").Value = 52000010
$ws.Range("N18").Value = -52000354
$ws.Range("H61").Value = 2084.8948
$ws.Range("I61").Value = 2096.8333
$ws.Range("K61").Value = 2096.8333
$ws.Range("M61").Value = -1894.8333
$ws.Range("H93").Value = 1416.0741
$ws.Range("I93").Value = 1302.2632
$ws.Range("J93").Value = 1686.375
$ws.Range("K93").Value = 1302.2632
$ws.Range("L93").Value = 1686.375
$ws.Range("M93").Value = -54.2632000000001
$ws.Range("N93").Value = -4182.375
$ws.Range("H113").Value = 2084.8948
$ws.Range("I113").Value = 2096.8333
$ws.Range("K113").Value = 2096.8333
$ws.Range("M113").Value = 73.16670000000022
$ws.Range("H122").Value = 4074.7097
$ws.Range("J122").Value = 3561.8
$ws.Range("L122").Value = 10685.4
$ws.Range("N122").Value = -15585.4

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H52").Value = 40000
$ws.Range("I52").Value = 0
$ws.Range("J52").Value = 40000
$ws.Range("K52").Value = 0
$ws.Range("L52").Value = 40000
$ws.Range("M52").ClearContents()
$ws.Range("N52").Value = -40452
$ws.Range("H54").Value = 38897.215
$ws.Range("I54").Value = 38445.25
$ws.Range("J54").Value = 39499.832
$ws.Range("K54").Value = 38445.25
$ws.Range("L54").Value = 39499.832
$ws.Range("M54").Value = -37925.25
$ws.Range("N54").Value = -40539.832
